$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Through 2022-07-09"

# Update the header label (shared string "2022 (through 07-08)" -> "2022 (through 07-09)")
$ws.Range("I1").Value = "2022 (through 07-09)"

# Update data cells
$ws.Range("I8").Value = 40
$ws.Range("I14").Value = 846
